$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) — update "want to go" counts (column F)
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(4, 6).Value = 866    # was 864
$ws.Cells.Item(5, 6).Value = 894    # was 891
$ws.Cells.Item(6, 6).Value = 1628   # was 1626
$ws.Cells.Item(7, 6).Value = 344    # was 342
$ws.Cells.Item(8, 6).Value = 1101   # was 1099
$ws.Cells.Item(9, 6).Value = 40     # was 39
$ws.Cells.Item(14, 6).Value = 101   # was 100
$ws.Cells.Item(15, 6).Value = 61    # was 60
$ws.Cells.Item(19, 6).Value = 37    # was 35
$ws.Cells.Item(20, 6).Value = 613   # was 612
$ws.Cells.Item(21, 6).Value = 604   # was 603
$ws.Cells.Item(23, 6).Value = 23    # was 22
$ws.Cells.Item(28, 6).Value = 2     # was 1

# Sheet "演出" (performances) — update "want to go" counts (column F)
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(5, 6).Value = 218    # was 217
$ws.Cells.Item(8, 6).Value = 103    # was 102
$ws.Cells.Item(9, 6).Value = 16     # was 14

# Sheet "全部类型" (all types) — update "want to go" counts (column F)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(5, 6).Value = 866    # was 864
$ws.Cells.Item(6, 6).Value = 894    # was 891
$ws.Cells.Item(7, 6).Value = 1628   # was 1626
$ws.Cells.Item(8, 6).Value = 344    # was 342
$ws.Cells.Item(9, 6).Value = 1101   # was 1099
$ws.Cells.Item(10, 6).Value = 40    # was 39
$ws.Cells.Item(15, 6).Value = 101   # was 100
$ws.Cells.Item(16, 6).Value = 61    # was 60
$ws.Cells.Item(23, 6).Value = 218   # was 217
$ws.Cells.Item(24, 6).Value = 218   # was 217
$ws.Cells.Item(25, 6).Value = 37    # was 35
$ws.Cells.Item(26, 6).Value = 613   # was 612
$ws.Cells.Item(27, 6).Value = 604   # was 603
$ws.Cells.Item(29, 6).Value = 23    # was 22
$ws.Cells.Item(36, 6).Value = 103   # was 102
$ws.Cells.Item(37, 6).Value = 103   # was 102
$ws.Cells.Item(38, 6).Value = 2     # was 1
$ws.Cells.Item(39, 6).Value = 16    # was 14
